$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 999.8
$ws.Cells.Item(127, 9).Value = 999.8
$ws.Cells.Item(127, 11).Value = 2999.4
$ws.Cells.Item(127, 13).Value = 1960.6

$ws.Cells.Item(138, 8).Value = 1405.7794
$ws.Cells.Item(138, 9).Value = 781.03125
$ws.Cells.Item(138, 10).Value = 1961.1111
$ws.Cells.Item(138, 11).Value = 2343.09375
$ws.Cells.Item(138, 12).Value = 5883.3333
$ws.Cells.Item(138, 13).Value = 2796.90625
$ws.Cells.Item(138, 14).Value = -16163.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3122.2976
$ws.Cells.Item(32, 9).Value = 2739.1948
$ws.Cells.Item(32, 10).Value = 7336.4287
$ws.Cells.Item(32, 11).Value = 2739.1948
$ws.Cells.Item(32, 12).Value = 7336.4287
$ws.Cells.Item(32, 13).Value = -2452.1948
$ws.Cells.Item(32, 14).Value = -7910.4287

$ws.Cells.Item(39, 8).Value = 2838.3333
$ws.Cells.Item(39, 9).Value = 2838.3333
$ws.Cells.Item(39, 11).Value = 2838.3333
$ws.Cells.Item(39, 13).Value = -2318.3333

$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).ClearContents()

$ws.Cells.Item(88, 8).Value = 1009801.9
$ws.Cells.Item(88, 9).Value = 1673668.6
$ws.Cells.Item(88, 10).Value = 14001.75
$ws.Cells.Item(88, 11).Value = 1673668.6
$ws.Cells.Item(88, 12).Value = 14001.75
$ws.Cells.Item(88, 13).Value = -1673262.6
$ws.Cells.Item(88, 14).Value = -14813.75

$ws.Cells.Item(91, 8).Value = 1009801.9
$ws.Cells.Item(91, 9).Value = 1673668.6
$ws.Cells.Item(91, 10).Value = 14001.75
$ws.Cells.Item(91, 11).Value = 1673668.6
$ws.Cells.Item(91, 12).Value = 14001.75
$ws.Cells.Item(91, 13).Value = -1672264.6
$ws.Cells.Item(91, 14).Value = -16809.75

$ws.Cells.Item(122, 8).Value = 942.1579
$ws.Cells.Item(122, 9).Value = 868
$ws.Cells.Item(122, 10).Value = 1008.9
$ws.Cells.Item(122, 11).Value = 2604
$ws.Cells.Item(122, 12).Value = 3026.7
$ws.Cells.Item(122, 13).Value = -154
$ws.Cells.Item(122, 14).Value = -7926.7

$ws.Cells.Item(132, 8).Value = 2848.0645
$ws.Cells.Item(132, 9).Value = 3264.1191
$ws.Cells.Item(132, 10).Value = 1974.35
$ws.Cells.Item(132, 11).Value = 9792.3573
$ws.Cells.Item(132, 12).Value = 5923.049999999999
$ws.Cells.Item(132, 13).Value = -7262.3573
$ws.Cells.Item(132, 14).Value = -10983.05

$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 13).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2843.3928
$ws.Cells.Item(105, 9).Value = 2947.7778
$ws.Cells.Item(105, 11).Value = 2947.7778
$ws.Cells.Item(105, 13).Value = -1200.7778

$ws.Cells.Item(107, 8).Value = 2230
$ws.Cells.Item(107, 9).Value = 2268.3333
$ws.Cells.Item(107, 11).Value = 2268.3333
$ws.Cells.Item(107, 13).Value = -348.3332999999998

$ws.Cells.Item(134, 8).Value = 25219.256
$ws.Cells.Item(134, 9).Value = 37407.855
$ws.Cells.Item(134, 10).Value = 2467.2
$ws.Cells.Item(134, 11).Value = 112223.565
$ws.Cells.Item(134, 12).Value = 7401.599999999999
$ws.Cells.Item(134, 13).Value = -109688.565
$ws.Cells.Item(134, 14).Value = -12471.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6413328
$ws.Cells.Item(31, 9).Value = 2607.6155
$ws.Cells.Item(31, 10).Value = 12824048
$ws.Cells.Item(31, 11).Value = 2607.6155
$ws.Cells.Item(31, 12).Value = 12824048
$ws.Cells.Item(31, 13).Value = -2312.6155
$ws.Cells.Item(31, 14).Value = -12824638

$ws.Cells.Item(34, 8).Value = 6413328
$ws.Cells.Item(34, 9).Value = 2607.6155
$ws.Cells.Item(34, 10).Value = 12824048
$ws.Cells.Item(34, 11).Value = 2607.6155
$ws.Cells.Item(34, 12).Value = 12824048
$ws.Cells.Item(34, 13).Value = -2405.6155
$ws.Cells.Item(34, 14).Value = -12824452

$ws.Cells.Item(62, 8).Value = 55557732
$ws.Cells.Item(62, 9).Value = 2151.6667
$ws.Cells.Item(62, 10).Value = 111113310
$ws.Cells.Item(62, 11).Value = 2151.6667
$ws.Cells.Item(62, 12).Value = 111113310
$ws.Cells.Item(62, 13).Value = -1527.6667
$ws.Cells.Item(62, 14).Value = -111114558

$ws.Cells.Item(65, 8).Value = 55557732
$ws.Cells.Item(65, 9).Value = 2151.6667
$ws.Cells.Item(65, 10).Value = 111113310
$ws.Cells.Item(65, 11).Value = 10758.3335
$ws.Cells.Item(65, 12).Value = 555566550
$ws.Cells.Item(65, 13).Value = -7638.333500000001
$ws.Cells.Item(65, 14).Value = -555572790

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 10451.363
$ws.Cells.Item(99, 10).Value = 18333.166
$ws.Cells.Item(99, 12).Value = 54999.49800000001
$ws.Cells.Item(99, 14).Value = -59491.49800000001

$ws.Cells.Item(122, 8).Value = 989.2593
$ws.Cells.Item(122, 9).Value = 1232.6666
$ws.Cells.Item(122, 10).Value = 794.5333
$ws.Cells.Item(122, 11).Value = 11093.9994
$ws.Cells.Item(122, 12).Value = 7150.7997
$ws.Cells.Item(122, 13).Value = -8643.9994
$ws.Cells.Item(122, 14).Value = -12050.7997

$ws.Cells.Item(125, 8).Value = 4550

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 24289470
$ws.Cells.Item(70, 9).Value = 28337076
$ws.Cells.Item(70, 10).Value = 3833.3333
$ws.Cells.Item(70, 11).Value = 28337076
$ws.Cells.Item(70, 12).Value = 3833.3333
$ws.Cells.Item(70, 13).Value = -28336806
$ws.Cells.Item(70, 14).Value = -4373.3333

$ws.Cells.Item(73, 8).Value = 24289470
$ws.Cells.Item(73, 9).Value = 28337076
$ws.Cells.Item(73, 10).Value = 3833.3333
$ws.Cells.Item(73, 11).Value = 28337076
$ws.Cells.Item(73, 12).Value = 3833.3333
$ws.Cells.Item(73, 13).Value = -28336140
$ws.Cells.Item(73, 14).Value = -5705.3333

$ws.Cells.Item(102, 8).Value = 1223
$ws.Cells.Item(102, 9).Value = 1029.5714
$ws.Cells.Item(102, 11).Value = 1029.5714
$ws.Cells.Item(102, 13).Value = 592.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2558.3333
$ws.Cells.Item(7, 9).Value = 1233.3334
$ws.Cells.Item(7, 10).Value = 3000
$ws.Cells.Item(7, 11).Value = 1233.3334
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = -1121.3334
$ws.Cells.Item(7, 14).Value = -3224

$ws.Cells.Item(126, 8).Value = 2558.3333
$ws.Cells.Item(126, 9).Value = 1233.3334
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 3700.0002
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -1230.0002
$ws.Cells.Item(126, 14).Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1000
$ws.Cells.Item(96, 9).Value = 1000
$ws.Cells.Item(96, 11).Value = 1000
$ws.Cells.Item(96, 13).Value = 373

$ws.Cells.Item(132, 8).Value = 4160.3
$ws.Cells.Item(132, 9).Value = 3950
$ws.Cells.Item(132, 10).Value = 4300.5
$ws.Cells.Item(132, 11).Value = 11850
$ws.Cells.Item(132, 12).Value = 12901.5
$ws.Cells.Item(132, 13).Value = -9320
$ws.Cells.Item(132, 14).Value = -17961.5
